# Generate Report for Handoff
#
# Re-generating the handoff report updates the "Latest Handoff Datetime"
# (column D) for the file that was just handed off
# (07b190c3-b3e5-4230-a29e-3ad6e0bd6d5e.md) on both locale sheets.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D4").Value = "2016-01-27 02:12:13"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D4").Value = "2016-01-27 02:12:24"
